# Notes from the meeting today with SOM:
# Sort the VA_waittimes data (rows 2:37, columns A:R) alphabetically by the
# "Location" column (C), ascending, and leave the resulting range selected
# the way Excel does right after running Data > Sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VA_waittimes")
$ws.Activate()

$sortRange = $ws.Range("A2:R37")
$sortKey = $ws.Range("C2")

$sortRange.Sort($sortKey, 1, $null, $null, 2)

$ws.Range("A2:R37").Select() | Out-Null
